$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: coin name + link change (plain text, safe to set directly)
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

# Price (D) and Volume(1h) (E) columns: values look numeric/date-like to Excel's
# smart-entry parser, so we stage them via a text formula and then convert the
# formula result to a static value with PasteSpecial (values only) -- this keeps
# the cells as plain text, matching the original inline-string content, and never
# touches any NumberFormat/style.

$ws.Range("D2").Formula = '="29.288.76"'
$ws.Range("E2").Formula = '="  +0.36%  "'
$ws.Range("D3").Formula = '="1.871.48"'
$ws.Range("E3").Formula = '="  +0.27%  "'
$ws.Range("E4").Formula = '="  +0.16%  "'
$ws.Range("D5").Formula = '="0.7061"'
$ws.Range("E5").Formula = '="  -0.75%  "'
$ws.Range("D6").Formula = '="241.42"'
$ws.Range("E6").Formula = '="  -0.05%  "'
$ws.Range("D7").Formula = '="1.000"'
$ws.Range("E7").Formula = '="  +0.05%  "'
$ws.Range("D8").Formula = '="0.3098"'
$ws.Range("E8").Formula = '="  -0.17%  "'
$ws.Range("D9").Formula = '="0.07746"'
$ws.Range("E9").Formula = '="  +1.13%  "'
$ws.Range("D10").Formula = '="25.07"'
$ws.Range("E10").Formula = '="  +1.40%  "'
$ws.Range("D11").Formula = '="0.08374"'
$ws.Range("E11").Formula = '="  +0.23%  "'
$ws.Range("D12").Formula = '="1.869.41"'
$ws.Range("E12").Formula = '="  -0.09%  "'
$ws.Range("D13").Formula = '="5.234"'
$ws.Range("E13").Formula = '="  +0.02%  "'
$ws.Range("D14").Formula = '="0.7103"'
$ws.Range("E14").Formula = '="  +0.05%  "'
$ws.Range("D15").Formula = '="91.01"'
$ws.Range("E15").Formula = '="  -0.37%  "'
$ws.Range("D16").Formula = '="29.300.27"'
$ws.Range("E16").Formula = '="  +0.32%  "'
$ws.Range("D17").Formula = '="6.047"'
$ws.Range("E17").Formula = '="  +1.92%  "'
$ws.Range("D18").Formula = '="0.000008157"'
$ws.Range("E18").Formula = '="  +4.35%  "'
$ws.Range("D19").Formula = '="239.63"'
$ws.Range("E19").Formula = '="  -1.67%  "'
$ws.Range("E20").Formula = '="  +0.70%  "'
$ws.Range("D21").Formula = '="2.119.21"'
$ws.Range("E21").Formula = '="  +0.13%  "'
$ws.Range("E22").Formula = '="  +0.22%  "'
$ws.Range("D23").Formula = '="7.735"'
$ws.Range("E23").Formula = '="  -1.71%  "'
$ws.Range("D24").Formula = '="1.001"'
$ws.Range("E24").Formula = '="  +0.16%  "'
$ws.Range("E25").Formula = '="  -0.54%  "'
$ws.Range("D26").Formula = '="162.81"'
$ws.Range("E26").Formula = '="  -0.33%  "'
$ws.Range("D27").Formula = '="9.006"'
$ws.Range("E27").Formula = '="  +0.64%  "'
$ws.Range("D28").Formula = '="18.48"'
$ws.Range("E28").Formula = '="  +0.00%  "'
$ws.Range("D29").Formula = '="1.508"'
$ws.Range("E29").Formula = '="  +0.59%  "'
$ws.Range("D30").Formula = '="4.397"'
$ws.Range("E30").Formula = '="  -0.05%  "'
$ws.Range("D31").Formula = '="1.291"'
$ws.Range("E31").Formula = '="  -2.46%  "'
$ws.Range("D32").Formula = '="4.299"'
$ws.Range("E32").Formula = '="  +1.29%  "'
$ws.Range("D33").Formula = '="0.05290"'
$ws.Range("E33").Formula = '="  +2.52%  "'
$ws.Range("D34").Formula = '="1.935"'
$ws.Range("E34").Formula = '="  +1.07%  "'
$ws.Range("D35").Formula = '="1.176"'
$ws.Range("E35").Formula = '="  +0.83%  "'
$ws.Range("D36").Formula = '="0.7418"'
$ws.Range("E36").Formula = '="  -7.60%  "'
$ws.Range("E37").Formula = '="  +0.81%  "'
$ws.Range("D38").Formula = '="0.01869"'
$ws.Range("E38").Formula = '="  +1.10%  "'
$ws.Range("D39").Formula = '="1.226.34"'
$ws.Range("E39").Formula = '="  +5.27%  "'
$ws.Range("D40").Formula = '="2.730"'
$ws.Range("E40").Formula = '="  +0.83%  "'
$ws.Range("D41").Formula = '="6.547"'
$ws.Range("E41").Formula = '="  +4.54%  "'
$ws.Range("D42").Formula = '="0.8838"'
$ws.Range("E42").Formula = '="  -0.85%  "'
$ws.Range("D43").Formula = '="72.35"'
$ws.Range("E43").Formula = '="  -1.05%  "'
$ws.Range("D44").Formula = '="108.90"'
$ws.Range("E44").Formula = '="  +6.02%  "'
$ws.Range("E45").Formula = '="  +0.12%  "'
$ws.Range("D46").Formula = '="2.016.09"'
$ws.Range("E46").Formula = '="  +0.27%  "'
$ws.Range("D47").Formula = '="0.5195"'
$ws.Range("E47").Formula = '="  +0.30%  "'
$ws.Range("D48").Formula = '="1.790"'
$ws.Range("E48").Formula = '="  +0.75%  "'
$ws.Range("E49").Formula = '="  +2.27%  "'
$ws.Range("D50").Formula = '="9.382"'
$ws.Range("E50").Formula = '="  +0.54%  "'
$ws.Range("D51").Formula = '="0.4301"'
$ws.Range("E51").Formula = '="  +0.34%  "'

# Convert the staged formulas to static text values in one batch operation.
$ws.Range("D2:E51").Copy()
$ws.Range("D2:E51").PasteSpecial(-4163)
$excel.CutCopyMode = 0

